# Scheduled market-data refresh: update cached price/profit figures (columns H:N)
# across the Leve-profit sheets. Values sourced from the latest market snapshot;
# a handful of rows also gain/lose the LeveProfitNQ/HQ (M/N) cells entirely.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 214.78572
$ws.Range("I9").Value = 146.41667
$ws.Range("K9").Value = 146.41667
$ws.Range("M9").Value = 22.58332999999999
$ws.Range("H32").Value = 3566.3333
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()
$ws.Range("H33").Value = 385.46155
$ws.Range("I33").Value = 194.78947
$ws.Range("K33").Value = 194.78947
$ws.Range("M33").Value = 34.21053000000001
$ws.Range("H41").Value = 584.6
$ws.Range("I41").Value = 66.42856999999999
$ws.Range("J41").Value = 1793.6666
$ws.Range("K41").Value = 66.42856999999999
$ws.Range("L41").Value = 1793.6666
$ws.Range("M41").Value = 373.57143
$ws.Range("N41").Value = -2673.6666
$ws.Range("H64").Value = 7711.316
$ws.Range("I64").Value = 5387
$ws.Range("K64").Value = 5387
$ws.Range("M64").Value = -5139
$ws.Range("H67").Value = 7711.316
$ws.Range("I67").Value = 5387
$ws.Range("K67").Value = 5387
$ws.Range("M67").Value = -4529
$ws.Range("H86").Value = 2218.9167
$ws.Range("I86").Value = 1135.3334
$ws.Range("J86").Value = 3302.5
$ws.Range("K86").Value = 1135.3334
$ws.Range("L86").Value = 3302.5
$ws.Range("M86").Value = -12.33339999999998
$ws.Range("N86").Value = -5548.5
$ws.Range("H89").Value = 2218.9167
$ws.Range("I89").Value = 1135.3334
$ws.Range("J89").Value = 3302.5
$ws.Range("K89").Value = 5676.666999999999
$ws.Range("L89").Value = 16512.5
$ws.Range("M89").Value = -60.66699999999946
$ws.Range("N89").Value = -27744.5
$ws.Range("H99").Value = 1952.6364
$ws.Range("I99").Value = 275
$ws.Range("K99").Value = 825
$ws.Range("M99").Value = 673
$ws.Range("H103").Value = 1313.6666
$ws.Range("I103").Value = 970.5
$ws.Range("J103").Value = 2000
$ws.Range("K103").Value = 2911.5
$ws.Range("L103").Value = 6000
$ws.Range("M103").Value = -2325.5
$ws.Range("N103").Value = -7172
$ws.Range("H129").Value = 1793.75
$ws.Range("I129").Value = 1478.5714
$ws.Range("J129").Value = 4000
$ws.Range("K129").Value = 4435.7142
$ws.Range("L129").Value = 12000
$ws.Range("M129").Value = 564.2857999999997
$ws.Range("N129").Value = -22000
$ws.Range("H132").Value = 53727.71
$ws.Range("I132").Value = 58345.465
$ws.Range("K132").Value = 175036.395
$ws.Range("M132").Value = -172506.395
$ws.Range("H135").Value = 1134.0646
$ws.Range("J135").Value = 1457.7858
$ws.Range("L135").Value = 13120.0722
$ws.Range("N135").Value = -18190.0722
$ws.Range("H137").Value = 1414394.6
$ws.Range("J137").Value = 3250074.2
$ws.Range("L137").Value = 9750222.600000001
$ws.Range("N137").Value = -9755322.600000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2532.3333
$ws.Range("I122").Value = 2532.3333
$ws.Range("K122").Value = 7596.999899999999
$ws.Range("M122").Value = -5146.999899999999
$ws.Range("H132").Value = 427751.1
$ws.Range("I132").Value = 514672.56
$ws.Range("K132").Value = 1544017.68
$ws.Range("M132").Value = -1541487.68

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H15").Value = 103
$ws.Range("I15").Value = 6
$ws.Range("J15").Value = 200
$ws.Range("K15").Value = 6
$ws.Range("L15").Value = 200
$ws.Range("M15").Value = 221
$ws.Range("N15").Value = -654
$ws.Range("H86").Value = 1456.8
$ws.Range("I86").Value = 1492.75
$ws.Range("J86").Value = 953.5
$ws.Range("K86").Value = 1492.75
$ws.Range("L86").Value = 953.5
$ws.Range("M86").Value = -369.75
$ws.Range("N86").Value = -3199.5
$ws.Range("H89").Value = 1456.8
$ws.Range("I89").Value = 1492.75
$ws.Range("J89").Value = 953.5
$ws.Range("K89").Value = 7463.75
$ws.Range("L89").Value = 4767.5
$ws.Range("M89").Value = -1847.75
$ws.Range("N89").Value = -15999.5
$ws.Range("H99").Value = 3967
$ws.Range("I99").Value = 4562.2
$ws.Range("J99").Value = 2975
$ws.Range("K99").Value = 4562.2
$ws.Range("L99").Value = 2975
$ws.Range("M99").Value = -3064.2
$ws.Range("N99").Value = -5971
$ws.Range("H105").Value = 8669
$ws.Range("I105").Value = 8669
$ws.Range("K105").Value = 8669
$ws.Range("M105").Value = -6922
$ws.Range("H107").Value = 4274.636
$ws.Range("I107").Value = 3773.1667
$ws.Range("J107").Value = 4876.4
$ws.Range("K107").Value = 3773.1667
$ws.Range("L107").Value = 4876.4
$ws.Range("M107").Value = -1853.1667
$ws.Range("N107").Value = -8716.4
$ws.Range("H134").Value = 1185964
$ws.Range("I134").Value = 1325871
$ws.Range("J134").Value = 766242.7
$ws.Range("K134").Value = 3977613
$ws.Range("L134").Value = 2298728.1
$ws.Range("M134").Value = -3975078
$ws.Range("N134").Value = -2303798.1

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 728.6667
$ws.Range("J16").Value = 1049
$ws.Range("L16").Value = 1049
$ws.Range("N16").Value = -1623
$ws.Range("H113").Value = 728.6667
$ws.Range("J113").Value = 1049
$ws.Range("L113").Value = 1049
$ws.Range("N113").Value = -5389
$ws.Range("H122").Value = 3424.8708
$ws.Range("J122").Value = 3679.3572
$ws.Range("L122").Value = 11038.0716
$ws.Range("N122").Value = -15938.0716

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 85047510
$ws.Range("I4").Value = 85047510
$ws.Range("K4").Value = 255142530
$ws.Range("M4").Value = -255142418
$ws.Range("H56").Value = 7659.5
$ws.Range("I56").Value = 7659.5
$ws.Range("K56").Value = 7659.5
$ws.Range("M56").Value = -7129.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 793230.75
$ws.Range("I21").Value = 10012000
$ws.Range("K21").Value = 10012000
$ws.Range("M21").Value = -10011827
$ws.Range("H30").Value = 793230.75
$ws.Range("I30").Value = 10012000
$ws.Range("K30").Value = 10012000
$ws.Range("M30").Value = -10011895
$ws.Range("H70").Value = 4769.7144
$ws.Range("I70").Value = 4797.6
$ws.Range("K70").Value = 4797.6
$ws.Range("M70").Value = -4527.6
$ws.Range("H73").Value = 4769.7144
$ws.Range("I73").Value = 4797.6
$ws.Range("K73").Value = 4797.6
$ws.Range("M73").Value = -3861.6
$ws.Range("H113").Value = 7434.875
$ws.Range("I113").Value = 4079.8333
$ws.Range("K113").Value = 4079.8333
$ws.Range("M113").Value = -1909.8333
$ws.Range("H126").Value = 879520.9399999999
$ws.Range("I126").Value = 1390697.6
$ws.Range("J126").Value = 3218.1428
$ws.Range("K126").Value = 4172092.8
$ws.Range("L126").Value = 9654.428400000001
$ws.Range("M126").Value = -4169622.8
$ws.Range("N126").Value = -14594.4284
$ws.Range("H132").Value = 25954778
$ws.Range("I132").Value = 30671956
$ws.Range("J132").Value = 10301.833
$ws.Range("K132").Value = 92015868
$ws.Range("L132").Value = 30905.499
$ws.Range("M132").Value = -92013338
$ws.Range("N132").Value = -35965.499

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1826.5
$ws.Range("I22").Value = 980
$ws.Range("J22").Value = 2249.75
$ws.Range("K22").Value = 980
$ws.Range("L22").Value = 2249.75
$ws.Range("M22").Value = -685
$ws.Range("N22").Value = -2839.75
$ws.Range("H23").Value = 16599
$ws.Range("I23").Value = 4898.5
$ws.Range("K23").Value = 4898.5
$ws.Range("M23").Value = -4668.5
$ws.Range("H27").Value = 1826.5
$ws.Range("I27").Value = 980
$ws.Range("J27").Value = 2249.75
$ws.Range("K27").Value = 980
$ws.Range("L27").Value = 2249.75
$ws.Range("M27").Value = -873
$ws.Range("N27").Value = -2463.75
$ws.Range("H46").Value = 6964.5835
$ws.Range("I46").Value = 9547.25
$ws.Range("J46").Value = 1799.25
$ws.Range("K46").Value = 9547.25
$ws.Range("L46").Value = 1799.25
$ws.Range("M46").Value = -9359.25
$ws.Range("N46").Value = -2175.25
$ws.Range("H61").Value = 2324.5
$ws.Range("I61").Value = 1560.8
$ws.Range("J61").Value = 6143
$ws.Range("K61").Value = 1560.8
$ws.Range("L61").Value = 6143
$ws.Range("M61").Value = -1358.8
$ws.Range("N61").Value = -6547
$ws.Range("H100").Value = 8147.353
$ws.Range("I100").Value = 2100.0833
$ws.Range("J100").Value = 22660.8
$ws.Range("K100").Value = 2100.0833
$ws.Range("L100").Value = 22660.8
$ws.Range("M100").Value = -1559.0833
$ws.Range("N100").Value = -23742.8
$ws.Range("H113").Value = 2324.5
$ws.Range("I113").Value = 1560.8
$ws.Range("J113").Value = 6143
$ws.Range("K113").Value = 1560.8
$ws.Range("L113").Value = 6143
$ws.Range("M113").Value = 609.2
$ws.Range("N113").Value = -10483

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 75000
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 75000
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 75000
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -75226
$ws.Range("H122").Value = 2782.742
$ws.Range("I122").Value = 2399.0386
$ws.Range("K122").Value = 7197.1158
$ws.Range("M122").Value = -4747.1158
$ws.Range("H126").Value = 5071.6113
$ws.Range("I126").Value = 4752.6
$ws.Range("J126").Value = 6666.6665
$ws.Range("K126").Value = 14257.8
$ws.Range("L126").Value = 19999.9995
$ws.Range("M126").Value = -11787.8
$ws.Range("N126").Value = -24939.9995
$ws.Range("H136").Value = 7385.4863
$ws.Range("I136").Value = 7567.375
$ws.Range("J136").Value = 6221.4
$ws.Range("K136").Value = 22702.125
$ws.Range("L136").Value = 18664.2
$ws.Range("M136").Value = -20152.125
$ws.Range("N136").Value = -23764.2
